# Update the "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptocurrency list on the active worksheet, per the scheduled GitHub
# Actions data refresh. Values are entered with a leading apostrophe so
# Excel stores them as literal text (matching the original inline-string
# cells) instead of auto-converting look-alike numbers (e.g. "308.60" or
# "1.010") and losing significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'26.554.16"
$ws.Cells.Item(2, 5).Value = "'  -3.08%  "
$ws.Cells.Item(3, 4).Value = "'1.805.31"
$ws.Cells.Item(3, 5).Value = "'  -2.84%  "
$ws.Cells.Item(4, 4).Value = "'1.011"
$ws.Cells.Item(4, 5).Value = "'  +0.63%  "
$ws.Cells.Item(5, 4).Value = "'1.010"
$ws.Cells.Item(5, 5).Value = "'  +0.65%  "
$ws.Cells.Item(6, 4).Value = "'308.60"
$ws.Cells.Item(6, 5).Value = "'  -1.95%  "
$ws.Cells.Item(7, 4).Value = "'0.4542"
$ws.Cells.Item(7, 5).Value = "'  -1.90%  "
$ws.Cells.Item(8, 4).Value = "'0.3669"
$ws.Cells.Item(8, 5).Value = "'  -1.47%  "
$ws.Cells.Item(9, 4).Value = "'0.07132"
$ws.Cells.Item(9, 5).Value = "'  -2.56%  "
$ws.Cells.Item(10, 4).Value = "'0.8697"
$ws.Cells.Item(10, 5).Value = "'  -2.17%  "
$ws.Cells.Item(11, 4).Value = "'0.07783"
$ws.Cells.Item(11, 5).Value = "'  -0.78%  "
$ws.Cells.Item(12, 4).Value = "'19.22"
$ws.Cells.Item(12, 5).Value = "'  -3.80%  "
$ws.Cells.Item(13, 4).Value = "'1.837.49"
$ws.Cells.Item(13, 5).Value = "'  -2.30%  "
$ws.Cells.Item(14, 4).Value = "'5.283"
$ws.Cells.Item(14, 5).Value = "'  -2.31%  "
$ws.Cells.Item(15, 4).Value = "'6.319"
$ws.Cells.Item(15, 5).Value = "'  -3.89%  "
$ws.Cells.Item(16, 4).Value = "'86.67"
$ws.Cells.Item(16, 5).Value = "'  -5.72%  "
$ws.Cells.Item(17, 4).Value = "'1.012"
$ws.Cells.Item(17, 5).Value = "'  +0.72%  "
$ws.Cells.Item(18, 4).Value = "'0.000008578"
$ws.Cells.Item(18, 5).Value = "'  -4.50%  "
$ws.Cells.Item(19, 4).Value = "'1.010"
$ws.Cells.Item(19, 5).Value = "'  +0.58%  "
$ws.Cells.Item(20, 4).Value = "'26.586.16"
$ws.Cells.Item(20, 5).Value = "'  -3.04%  "
$ws.Cells.Item(21, 4).Value = "'14.24"
$ws.Cells.Item(21, 5).Value = "'  -3.74%  "
$ws.Cells.Item(22, 4).Value = "'4.958"
$ws.Cells.Item(22, 5).Value = "'  -3.39%  "
$ws.Cells.Item(23, 4).Value = "'2.081.06"
$ws.Cells.Item(23, 5).Value = "'  +0.69%  "
$ws.Cells.Item(24, 5).Value = "'  -1.89%  "
$ws.Cells.Item(25, 4).Value = "'1.986"
$ws.Cells.Item(25, 5).Value = "'  +2.75%  "
$ws.Cells.Item(26, 4).Value = "'151.05"
$ws.Cells.Item(26, 5).Value = "'  -0.63%  "
$ws.Cells.Item(27, 4).Value = "'17.95"
$ws.Cells.Item(27, 5).Value = "'  -2.78%  "
$ws.Cells.Item(28, 4).Value = "'1.976"
$ws.Cells.Item(28, 5).Value = "'  -4.12%  "
$ws.Cells.Item(29, 4).Value = "'113.02"
$ws.Cells.Item(29, 5).Value = "'  -2.74%  "
$ws.Cells.Item(30, 4).Value = "'4.870"
$ws.Cells.Item(30, 5).Value = "'  -4.54%  "
$ws.Cells.Item(31, 4).Value = "'0.08706"
$ws.Cells.Item(31, 5).Value = "'  -1.70%  "
$ws.Cells.Item(32, 4).Value = "'3.022"
$ws.Cells.Item(32, 5).Value = "'  -4.21%  "
$ws.Cells.Item(33, 4).Value = "'0.7322"
$ws.Cells.Item(33, 5).Value = "'  -4.96%  "
$ws.Cells.Item(34, 4).Value = "'4.438"
$ws.Cells.Item(34, 5).Value = "'  -1.76%  "
$ws.Cells.Item(35, 4).Value = "'1.111"
$ws.Cells.Item(35, 5).Value = "'  -5.65%  "
$ws.Cells.Item(36, 4).Value = "'2.477"
$ws.Cells.Item(36, 5).Value = "'  -7.88%  "
$ws.Cells.Item(37, 4).Value = "'1.078"
$ws.Cells.Item(37, 5).Value = "'  -0.35%  "
$ws.Cells.Item(38, 4).Value = "'0.01916"
$ws.Cells.Item(38, 5).Value = "'  -2.42%  "
$ws.Cells.Item(39, 4).Value = "'0.05091"
$ws.Cells.Item(39, 5).Value = "'  -2.88%  "
$ws.Cells.Item(40, 4).Value = "'2.863"
$ws.Cells.Item(40, 5).Value = "'  -3.49%  "
$ws.Cells.Item(41, 4).Value = "'6.877"
$ws.Cells.Item(41, 5).Value = "'  -2.98%  "
$ws.Cells.Item(42, 4).Value = "'0.4903"
$ws.Cells.Item(42, 5).Value = "'  -4.91%  "
$ws.Cells.Item(43, 4).Value = "'0.1568"
$ws.Cells.Item(43, 5).Value = "'  -4.42%  "
$ws.Cells.Item(44, 4).Value = "'8.106"
$ws.Cells.Item(44, 5).Value = "'  -3.75%  "
$ws.Cells.Item(45, 4).Value = "'1.011"
$ws.Cells.Item(45, 5).Value = "'  +0.77%  "
$ws.Cells.Item(46, 4).Value = "'0.4585"
$ws.Cells.Item(46, 5).Value = "'  -4.87%  "
$ws.Cells.Item(47, 4).Value = "'102.73"
$ws.Cells.Item(47, 5).Value = "'  -0.24%  "
$ws.Cells.Item(48, 4).Value = "'9.965"
$ws.Cells.Item(48, 5).Value = "'  -3.92%  "
$ws.Cells.Item(49, 4).Value = "'1.581"
$ws.Cells.Item(49, 5).Value = "'  -4.44%  "
$ws.Cells.Item(50, 4).Value = "'0.05994"
$ws.Cells.Item(50, 5).Value = "'  -3.60%  "
$ws.Cells.Item(51, 4).Value = "'63.54"
$ws.Cells.Item(51, 5).Value = "'  -2.64%  "
